# Roll the quarterly income-statement table forward by one quarter:
# drop the oldest quarter (column D) and shift every later quarter one
# column to the left, then fill in the newly-reported quarter in column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the column widths before we start moving things around, so the
# "wide" (published-date) column formatting can move left together with
# the data.
$widths = @{}
foreach ($c in @("E", "F", "G", "H", "I", "J", "K", "L", "M")) {
    $widths[$c] = $ws.Columns($c).ColumnWidth
}

# Shift the whole data block (period headers, publish dates and all the
# financial rows) one column to the left: E8:M27 -> D8:L27. PasteSpecial
# xlPasteAll carries the per-cell number formatting/styles along with the
# values, and naturally overwrites the old (oldest-quarter) column D data.
$ws.Range("E8:M27").Copy() | Out-Null
$ws.Range("D8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll) | Out-Null
$excel.CutCopyMode = 0

# Shift the column widths the same way so the "published date" column
# (the wider one) lines up with its data again.
$destCols = @("D", "E", "F", "G", "H", "I", "J", "K", "L")
$srcCols  = @("E", "F", "G", "H", "I", "J", "K", "L", "M")
for ($i = 0; $i -lt $destCols.Length; $i++) {
    $ws.Columns($destCols[$i]).ColumnWidth = $widths[$srcCols[$i]]
}
$ws.Columns("M").ColumnWidth = $widths["F"]

# New (latest) quarter header label and its publish date, now in column M.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30 (2)"

# The republished-statement footnote for an older quarter also moved:
# it now lands on column I after the shift, and its footnote count changed.
$ws.Range("I9").Value = "1402-02-30 (8)"

# New financial figures reported for the newest quarter (column M).
$ws.Range("M11").Value = 7406
$ws.Range("M12").Value = -5033
$ws.Range("M13").Value = 2372
$ws.Range("M14").Value = -792
$ws.Range("M16").Value = 1
$ws.Range("M17").Value = 1581
$ws.Range("M18").Value = -79
$ws.Range("M19").Value = 14
$ws.Range("M20").Value = 1516
$ws.Range("M21").Value = 83
$ws.Range("M22").Value = 1599
$ws.Range("M24").Value = 1599
$ws.Range("M26").Value = 9865
